$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$changes = @{
    "H86" = 45458490
    "J86" = 83338620
    "L86" = 83338620
    "N86" = -83340866
    "H89" = 45458490
    "J89" = 83338620
    "L89" = 416693100
    "N89" = -416704332
    "H106" = 2271.5625
    "I106" = 2103.4614
    "K106" = 2103.4614
    "M106" = -1472.4614
    "H107" = 588661.6
    "I107" = 714637.0600000001
    "K107" = 714637.0600000001
    "M107" = -712717.0600000001
    "H112" = 2250.8333
    "J112" = 2250.8333
    "L112" = 6752.499899999999
    "N112" = -8968.499899999999
    "H121" = 819.6087
    "J121" = 821.4761999999999
    "L121" = 2464.4286
    "N121" = -5958.428599999999
    "H129" = 1684557.4
    "I129" = 0
    "J129" = 1684557.4
    "K129" = 0
    "N129" = -5063672.199999999
    "H132" = 1738.2678
    "I132" = 1405.4
    "K132" = 4216.200000000001
    "M132" = -1686.200000000001
    "H137" = 7043973
    "I137" = 1353.1702
    "J137" = 20835770
    "K137" = 4059.5106
    "L137" = 62507310
    "M137" = -1509.5106
    "N137" = -62512410
    "H138" = 3014.65
    "I138" = 1333.8823
    "J138" = 3880.5
    "K138" = 4001.6469
    "L138" = 11641.5
    "M138" = 1138.3531
    "N138" = -21921.5
    "L129" = 5053672.199999999
}
foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}
$ws.Range("M129").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$changes = @{
    "H26" = 0
    "J26" = 0
    "H32" = 24885.479
    "I32" = 21846.777
    "K32" = 21846.777
    "M32" = -21559.777
    "H132" = 2145.05
    "I132" = 1536.6428
    "J132" = 3564.6667
    "K132" = 4609.928400000001
    "L132" = 10694.0001
    "M132" = -2079.928400000001
    "N132" = -15754.0001
    "L26" = 0
}
foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}
$ws.Range("N26").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$changes = @{
    "H26" = 8937
    "I26" = 8937
    "K26" = 8937
    "M26" = -8645
    "H96" = 9577
    "I96" = 2945.6667
    "K96" = 2945.6667
    "M96" = -199.6667000000002
    "H99" = 38463320
    "I99" = 52633224
    "J99" = 2151.7144
    "K99" = 52633224
    "L99" = 2151.7144
    "M99" = -52631726
    "N99" = -5147.7144
    "H105" = 15666.667
    "I105" = 8500
    "J105" = 30000
    "K105" = 8500
    "L105" = 30000
    "M105" = -6753
    "N105" = -33494
    "H107" = 11050.077
    "I107" = 882.25
    "J107" = 27318.6
    "K107" = 882.25
    "L107" = 27318.6
    "M107" = 1037.75
    "N107" = -31158.6
    "H134" = 61168.207
    "I134" = 2605.4
    "J134" = 223842.67
    "K134" = 7816.200000000001
    "L134" = 671528.01
    "M134" = -5281.200000000001
    "N134" = -676598.01
}
foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$changes = @{
    "H16" = 1139.1818
    "I16" = 876
    "J16" = 1599.75
    "K16" = 876
    "L16" = 1599.75
    "M16" = -589
    "N16" = -2173.75
    "H31" = 1762.4902
    "I31" = 943.26666
    "J31" = 2103.8333
    "K31" = 943.26666
    "L31" = 2103.8333
    "M31" = -648.26666
    "N31" = -2693.8333
    "H34" = 1762.4902
    "I34" = 943.26666
    "J34" = 2103.8333
    "K34" = 943.26666
    "L34" = 2103.8333
    "M34" = -741.26666
    "N34" = -2507.8333
    "H58" = 4753.1665
    "I58" = 1016.2308
    "J58" = 7610.8237
    "K58" = 1016.2308
    "L58" = 7610.8237
    "M58" = -813.2308
    "N58" = -8016.8237
    "H60" = 10000
    "J60" = 0
    "L60" = 0
    "H105" = 3300.9092
    "I105" = 2858.5715
    "J105" = 4075
    "K105" = 2858.5715
    "L105" = 4075
    "M105" = -1111.5715
    "N105" = -7569
    "H107" = 418.17392
    "I107" = 482.75
    "J107" = 404.57895
    "K107" = 482.75
    "L107" = 404.57895
    "M107" = 1437.25
    "N107" = -4244.57895
    "H113" = 1139.1818
    "I113" = 876
    "J113" = 1599.75
    "K113" = 876
    "L113" = 1599.75
    "M113" = 1294
    "N113" = -5939.75
    "H132" = 3325.8333
    "I132" = 1002
    "J132" = 5649.6665
    "K132" = 3006
    "L132" = 16948.9995
    "M132" = -476
    "N132" = -22008.9995
    "H136" = 4753.1665
    "I136" = 1016.2308
    "J136" = 7610.8237
    "K136" = 3048.6924
    "L136" = 22832.4711
    "M136" = -498.6923999999999
    "N136" = -27932.4711
    "H140" = 46140
    "J140" = 46140
    "L140" = 46140
    "N140" = -56500
}
foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}
$ws.Range("N60").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$changes = @{
    "H107" = 159917.33
    "I107" = 326.1905
    "J107" = 279610.7
    "K107" = 978.5715
    "L107" = 838832.1000000001
    "M107" = 941.4285
    "N107" = -842672.1000000001
    "H131" = 12527029
    "J131" = 2700.4285
    "L131" = 8101.2855
    "N131" = -18181.2855
}
foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$changes = @{
    "H102" = 954.375
    "I102" = 784.6667
    "J102" = 3500
    "K102" = 784.6667
    "L102" = 3500
    "M102" = 837.3333
    "N102" = -6744
    "H132" = 4095.818
    "I132" = 2408.8
    "J132" = 5501.6665
    "K132" = 7226.400000000001
    "L132" = 16504.9995
    "M132" = -4696.400000000001
    "N132" = -21564.9995
    "H138" = 37670.9
    "J138" = 37670.9
    "L138" = 37670.9
    "N138" = -47950.9
}
foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$changes = @{
    "H61" = 865.26666
    "I61" = 876.7143
    "K61" = 876.7143
    "M61" = -674.7143
    "H93" = 1709.1428
    "I93" = 980.5294
    "J93" = 2835.182
    "K93" = 980.5294
    "L93" = 2835.182
    "M93" = 267.4706
    "N93" = -5331.182
    "H100" = 10102837
    "I100" = 13890601
    "J100" = 2133.3333
    "K100" = 13890601
    "L100" = 2133.3333
    "M100" = -13890060
    "N100" = -3215.3333
    "H113" = 865.26666
    "I113" = 876.7143
    "K113" = 876.7143
    "M113" = 1293.2857
    "H132" = 5054431
    "I132" = 6064250.5
    "J132" = 5332.6665
    "K132" = 18192751.5
    "L132" = 15997.9995
    "M132" = -18190221.5
    "N132" = -21057.9995
    "H139" = 47705
    "J139" = 47705
    "L139" = 47705
    "N139" = -57985
}
foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$changes = @{
    "H113" = 1125.6
    "I113" = 1125.6
    "K113" = 3376.8
    "M113" = -1206.8
    "H132" = 2600.889
    "I132" = 1217.8334
    "J132" = 3292.4167
    "K132" = 3653.5002
    "L132" = 9877.250100000001
    "M132" = -1123.5002
    "N132" = -14937.2501
    "H136" = 2209.362
    "I136" = 2427.5144
    "J136" = 1877.3914
    "K136" = 7282.5432
    "L136" = 5632.174199999999
    "M136" = -4732.5432
    "N136" = -10732.1742
    "H138" = 46903
    "J138" = 46903
    "L138" = 46903
    "N138" = -57183
}
foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}

Write-Output "applied changes"